$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '255.53'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '3.97%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.13'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-5.43%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.206'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.59%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05860'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.93%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.705'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '0.80%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8708'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.63%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9570'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '11.78%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1412'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.28%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07123'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '0.55%'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.29%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09202'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.61%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001545'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.08%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0006092'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '1.76%'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.005809'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-3.04%'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.56%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.229'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.02%'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.41%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.63%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03448'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '3.01%'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.33%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.523'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '1.27%'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.62%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.65%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001221'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.03%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004563'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '9.42%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.12%'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001467'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '1.28%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03818'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '1.48%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005642'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.28%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1103'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '3.09%'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.01%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009839'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.97%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005415'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '2.62%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.11%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.09002'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.21%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002129'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-24.73%'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.11%'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.11%'
